# Add a new data row (row 20) to Sheet1, mirroring the existing Ciruela /
# Agrícola del Norte S.A. de Arica records, for the weekly Fruta/Hortaliza
# consolidation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 20

$ws.Cells.Item($row, 1).Value  = 1
$ws.Cells.Item($row, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item($row, 4).Value  = 44628
# Keep the same date number format used by the other rows in column D.
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
$ws.Cells.Item($row, 5).Value  = 15
$ws.Cells.Item($row, 6).Value  = "Fruta"
$ws.Cells.Item($row, 7).Value  = 100103
$ws.Cells.Item($row, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item($row, 9).Value  = 100103002
$ws.Cells.Item($row, 10).Value = "Ciruela"
$ws.Cells.Item($row, 11).Value = "Black Amber"
$ws.Cells.Item($row, 12).Value = "Segunda"
$ws.Cells.Item($row, 13).Value = 270
$ws.Cells.Item($row, 14).Value = 15000
$ws.Cells.Item($row, 15).Value = 16000
$ws.Cells.Item($row, 16).Value = 15500
$ws.Cells.Item($row, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 861
$ws.Cells.Item($row, 20).Value = 18
